$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set explicit custom widths on columns A and B (engine stores ColumnWidth + 5/6,
# snapped to the nearest 1/6 increment, matching Excel's pixel-width model -
# these inputs land on the target stored widths of 18.83203125 and 28).
$ws.Columns.Item(1).ColumnWidth = 17.998697916666668
$ws.Columns.Item(2).ColumnWidth = 27.166666666666668

# Move the selection to A9:B9 with A9 as the active cell
$ws.Range("A9:B9").Select()
